$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 620.92
$ws.Range("J33").Value = 911.6
$ws.Range("L33").Value = 911.6
$ws.Range("N33").Value = -1369.6
$ws.Range("H40").Value = 1559.4
$ws.Range("I40").Value = 1299.4286
$ws.Range("J40").Value = 2166
$ws.Range("K40").Value = 1299.4286
$ws.Range("L40").Value = 2166
$ws.Range("M40").Value = -1124.4286
$ws.Range("N40").Value = -2516
$ws.Range("H52").Value = 700
$ws.Range("I52").Value = 700
$ws.Range("K52").Value = 2100
$ws.Range("M52").Value = -1940
$ws.Range("H86").Value = 2391
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("H87").Value = 58731.25
$ws.Range("J87").Value = 97475
$ws.Range("L87").Value = 97475
$ws.Range("N87").Value = -99971
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H89").Value = 2391
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""
$ws.Range("H90").Value = 58731.25
$ws.Range("J90").Value = 97475
$ws.Range("L90").Value = 292425
$ws.Range("N90").Value = -304905
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""
$ws.Range("H100").Value = 2864.8333
$ws.Range("J100").Value = 3322.25
$ws.Range("L100").Value = 3322.25
$ws.Range("N100").Value = -4404.25
$ws.Range("H125").Value = 1666.3334
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
$ws.Range("H138").Value = 2988.4707
$ws.Range("I138").Value = 2864.6667
$ws.Range("K138").Value = 8594.000100000001
$ws.Range("M138").Value = -3454.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1942.5
$ws.Range("I2").Value = 1942.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1942.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1829.5
$ws.Range("N2").Value = ""
$ws.Range("H61").Value = 2308.8
$ws.Range("I61").Value = 2308.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2308.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2096.8
$ws.Range("N61").Value = ""
$ws.Range("H102").Value = 2997.8333
$ws.Range("I102").Value = 2995
$ws.Range("J102").Value = 2999.25
$ws.Range("K102").Value = 2995
$ws.Range("L102").Value = 2999.25
$ws.Range("M102").Value = -1373
$ws.Range("N102").Value = -6243.25
$ws.Range("H116").Value = 1942.5
$ws.Range("I116").Value = 1942.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1942.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 351.5
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 168647
$ws.Range("I132").Value = 201839.4
$ws.Range("K132").Value = 605518.2
$ws.Range("M132").Value = -602988.2
$ws.Range("H136").Value = 2308.8
$ws.Range("I136").Value = 2308.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6926.400000000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4376.400000000001
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1942.5
$ws.Range("I3").Value = 1942.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1942.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1828.5
$ws.Range("N3").Value = ""
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
$ws.Range("H88").Value = 29500
$ws.Range("J88").Value = 29500
$ws.Range("L88").Value = 29500
$ws.Range("N88").Value = -30312
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""
$ws.Range("H91").Value = 29500
$ws.Range("J91").Value = 29500
$ws.Range("L91").Value = 29500
$ws.Range("N91").Value = -32308
$ws.Range("H99").Value = 4837.077
$ws.Range("I99").Value = 4406.8335
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 4406.8335
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -2908.8335
$ws.Range("N99").Value = -12996
$ws.Range("H100").Value = 85000
$ws.Range("J100").Value = 85000
$ws.Range("L100").Value = 85000
$ws.Range("N100").Value = -87164
$ws.Range("H105").Value = 1983.1666
$ws.Range("I105").Value = 1983.1666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1983.1666
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -236.1666
$ws.Range("N105").Value = ""
$ws.Range("H134").Value = 7766.6665
$ws.Range("I134").Value = 7766.6665
$ws.Range("K134").Value = 23299.9995
$ws.Range("M134").Value = -20764.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1800
$ws.Range("H105").Value = 390
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 390
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 390
$ws.Range("M105").Value = ""
$ws.Range("N105").Value = -3884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 438
$ws.Range("J33").Value = 400
$ws.Range("L33").Value = 2400
$ws.Range("N33").Value = -2966
$ws.Range("H59").Value = 812.5
$ws.Range("I59").Value = 750
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 2250
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -1710
$ws.Range("N59").Value = -4080
$ws.Range("H75").Value = 5000
$ws.Range("J75").Value = 5000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 5000
$ws.Range("J78").Value = 5000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -54984
$ws.Range("H81").Value = 33250
$ws.Range("J81").Value = 37857.145
$ws.Range("L81").Value = 113571.435
$ws.Range("N81").Value = -115817.435
$ws.Range("H84").Value = 33250
$ws.Range("J84").Value = 37857.145
$ws.Range("L84").Value = 340714.305
$ws.Range("N84").Value = -351946.305
$ws.Range("H129").Value = 1830.1666
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 2097.2
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 6291.599999999999
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -16291.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 2008
$ws.Range("J19").Value = 2008
$ws.Range("L19").Value = 2008
$ws.Range("N19").Value = -2584
$ws.Range("H70").Value = 2998.5
$ws.Range("J70").Value = 2998
$ws.Range("L70").Value = 2998
$ws.Range("N70").Value = -3538
$ws.Range("H73").Value = 2998.5
$ws.Range("J73").Value = 2998
$ws.Range("L73").Value = 2998
$ws.Range("N73").Value = -4870
$ws.Range("H126").Value = 5399.5
$ws.Range("I126").Value = 5399.5
$ws.Range("K126").Value = 16198.5
$ws.Range("M126").Value = -13728.5
$ws.Range("H132").Value = 2698.3333
$ws.Range("I132").Value = 2801.75
$ws.Range("J132").Value = 2491.5
$ws.Range("K132").Value = 8405.25
$ws.Range("L132").Value = 7474.5
$ws.Range("M132").Value = -5875.25
$ws.Range("N132").Value = -12534.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""
$ws.Range("H20").Value = 4000
$ws.Range("J20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("N20").Value = -4452
$ws.Range("H22").Value = 843.125
$ws.Range("I22").Value = 618.9
$ws.Range("J22").Value = 1216.8334
$ws.Range("K22").Value = 618.9
$ws.Range("L22").Value = 1216.8334
$ws.Range("M22").Value = -323.9
$ws.Range("N22").Value = -1806.8334
$ws.Range("H27").Value = 843.125
$ws.Range("I27").Value = 618.9
$ws.Range("J27").Value = 1216.8334
$ws.Range("K27").Value = 618.9
$ws.Range("L27").Value = 1216.8334
$ws.Range("M27").Value = -511.9
$ws.Range("N27").Value = -1430.8334
$ws.Range("H46").Value = 598
$ws.Range("I46").Value = 598
$ws.Range("K46").Value = 598
$ws.Range("M46").Value = -410
$ws.Range("H68").Value = 4499.75
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 5999
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 5999
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -7497
$ws.Range("H71").Value = 4499.75
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 5999
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 29995
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -37483
$ws.Range("H82").Value = 3500
$ws.Range("I82").Value = 3500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3139
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 3500
$ws.Range("I85").Value = 3500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2252
$ws.Range("N85").Value = ""
